$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Valueprop")
[void]$ws.Activate()

# Update the admin test credentials (UserName / Password columns on row 2)
$ws.Range("B2").Value = "vraikanti@helenoftroy.com"
$ws.Range("C2").Value = "Baprvtyfaqtt2!"

# Move the active selection (matches the saved view state after the edit)
[void]$ws.Range("C8").Select()
